$d = $word.ActiveDocument

function Replace-Unique($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $old)
    }
}

# 1. Receipt timestamp at top of document. This run sits right next to a
# sibling run (" от ") that happens to share identical run formatting, so a
# plain text rewrite would let the two runs coalesce into one. Re-touching
# (no-op) a character property on just the replaced span keeps it as its
# own run, matching the original document's run layout.
$dateRng = $d.Content
$dateFound = $dateRng.Find.Execute("2017-03-15 03:57:05")
if (-not $dateFound) {
    Write-Output "NOT FOUND: 2017-03-15 03:57:05"
} else {
    $newDate = "2017-03-22 01:31:58"
    $dateStart = $dateRng.Start
    $dateRng.Text = $newDate
    $dateRng2 = $d.Range($dateStart, $dateStart + $newDate.Length)
    $origSize = $dateRng2.Font.Size
    $dateRng2.Font.Size = $origSize + 1
    $dateRng2.Font.Size = $origSize
}

# 2. Product article / description (English row)
Replace-Unique "Bolon, Create, 2x410x205" "Bolon Bolon by you  Individual order (over 600m2)  2000x25000 mm "

# 3. "Mest" (places) column
Replace-Unique "0|" "0.06|"

# 4. Packaging unit
Replace-Unique "package" "roll"

# 5. Quantity value (keep the "kg" run intact, only blank out the "2.9 " run)
Replace-Unique "2.9 " " "

# 6. Product article / description (Russian row)
Replace-Unique "Болон, Create, 2x410x205" "Виниловое покрытие Bolon Bolon by you (свыше 600 м2) 2000x25000 м"

# The remaining numeric cells ("0"->"3", "121"->"93", "0"->"279", "0"->"279")
# are too short/common for a plain text Find, so address them directly via
# the table/cell object model.
$t = $d.Tables.Item(2)

# Product row (row 4): m2 quantity, price, and line total.
$row = $t.Rows.Item(4)
$row.Cells.Item(10).Range.Text = "3"
$row.Cells.Item(12).Range.Text = "93"
$row.Cells.Item(13).Range.Text = "279"

# Totals row (row 7): grand total.
$totalsRow = $t.Rows.Item(7)
$totalsRow.Cells.Item(9).Range.Text = "279"
